$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") and column E ("Volume(1h)") values are updated row by
# row to match the refreshed cryptocurrency data. Column D entries that look
# like plain numbers are written with a leading apostrophe so Excel keeps
# storing them as text (matching the original inline-string cell type)
# instead of silently converting them to numeric values.

$ws.Range("D2").Value = "28.510.30"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.871.64"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -2.23%  "
$ws.Range("D5").Value = "'315.31"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("D7").Value = "'0.5082"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").Value = "'0.3903"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").Value = "'0.08356"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "'1.106"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").Value = "'41.84"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "'6.217"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").Value = "1.866.09"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "'20.43"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "'7.276"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").Value = "'0.00001100"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'0.06731"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "'1.006"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").Value = "'5.922"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").Value = "28.548.85"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'11.09"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "'2.195"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("D26").Value = "2.091.29"
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("D27").Value = "'158.74"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").Value = "'20.62"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("D29").Value = "'2.425"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("D30").Value = "'127.05"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("D32").Value = "'1.047"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "'5.735"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("D34").Value = "'3.613"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").Value = "'0.02462"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "'0.06602"
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").Value = "'0.2168"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").Value = "'8.894"
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("D39").Value = "'5.035"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "'1.181"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").Value = "'1.237"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "'0.6379"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").Value = "'11.10"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").Value = "'1.006"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").Value = "'0.6010"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").Value = "'13.14"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("D48").Value = "'2.011"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'122.57"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "'0.06809"
$ws.Range("E51").Value = "  -1.08%  "
